# Insert a new header row at the very top of Sheet1, pushing all existing
# data rows down by one (old row 1 becomes row 2, ..., old row 500 becomes
# row 501), then populate the new row 1 with the column headers "LABEL"
# and "TEXT".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a whole new row above row 1 - shifts everything down.
$ws.Rows("1:1").Insert()

# Set the new header row contents.
$ws.Range("A1").Value = "LABEL"
$ws.Range("B1").Value = "TEXT"

# Leave the selection on A2, matching the post-edit workbook state.
$ws.Range("A2").Select() | Out-Null
